$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the value in the last existing row (116), column B.
$ws.Cells.Item(116, 2).Value = 230.8219

# Append a new row (117) with the next monthly date label and its value.
$cellA117 = $ws.Cells.Item(117, 1)

# The label looks like a date ("01-08-2021"), so Excel would normally
# auto-convert it into a date serial value when typed directly. Temporarily
# mark the cell as Text so it is stored as a plain string (matching the
# existing date-label cells in column A), then restore the cell's style to
# the same (default/general) style used by the rest of the column so no
# visible formatting change is introduced.
$cellA117.NumberFormat = "@"
$cellA117.Value = "01-08-2021"
$cellA117.Style = $ws.Cells.Item(116, 1).Style

$ws.Cells.Item(117, 2).Value = 168.6209
